# "syllabus update for f22"
#
# 1. The "Well, guess what? There's also office hours..." sentence gets
#    split by Word's grammar checker: "There's" is wrapped in
#    proofErr (gramStart/gramEnd) markers, which forces that run to be
#    split into three runs.
# 2. The "Revised <date>" footer field result is bumped from
#    2020-08-23 to 2022-08-27.

$d = $word.ActiveDocument

# --- 1. Split the "There's" run and wrap it with grammar proofErr marks ---

$whole = "If you'd like to meet up to talk through something, I will always be happy to schedule a meeting via zoom or phone. That said, I've heard from many students that they feel weird about asking to take up my time (I hope you don't feel that way, I really do want to meet with you). Well, guess what? There's also office hours: Time each week already dedicated to you; no asking required!"

$target = $d.Content
$found = $target.Find.Execute($whole)

if ($found) {
    $start = $target.Start

    # Remove the whole paragraph's old single-run text...
    $target.Delete()

    # ...and drop in the same paragraph, now split into three runs with
    # proofErr markers bracketing "There's" (exactly what Word's grammar
    # checker does when it flags that word).
    $freshRange = $d.Range($start, $start)
    $xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000E0132" w:rsidRDefault="000E0132" w:rsidP="000E0132"><w:r><w:t xml:space="preserve">If you''d like to meet up to talk through something, I will always be happy to schedule a meeting via zoom or phone. That said, I''ve heard from many students that they feel weird about asking to take up my time (I hope you don''t feel that way, I really do want to meet with you). Well, guess what? </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>There''s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> also office hours: Time each week already dedicated to you; no asking required!</w:t></w:r></w:p>'
    $freshRange.InsertXML($xmlFrag)
}

# --- 2. Bump the "Revised" date shown in the default footer ---

$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("2020-08-23", $true, $false, $false, $false, $false, $true, 1, $false, "2022-08-27", 2) | Out-Null
